# Updated cryptos list on Mon Nov  4 07:17:14 UTC 2024 with GitHub Actions
#
# Price (column D) values are stored as text even though several of them
# look numeric ("560.41", "23.68", ...). Excel's COM layer auto-converts a
# plain numeric-looking string into a real number (losing the original
# formatting / trailing zeros, e.g. "1.00" -> 1). To keep those cells as
# plain text - matching the workbook's original inlineStr formatting - we
# prefix the assignment with a leading apostrophe (Excel's "force text"
# quote-prefix convention, stored but not part of the value) and then
# reset the cell style back to Normal so no stray style index is left on
# the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "69.049.53"
$ws.Range("E2").Value = "  +0.50%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.474.03"
$ws.Range("E3").Value = "  +0.74%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "560.41"
$ws.Range("E5").Value = "  -0.81%  "

# Row 6 - Solana
Set-TextValue "D6" "162.38"
$ws.Range("E6").Value = "  -0.97%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.06%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -2.73%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +1.26%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("E13").Value = "  -0.05%  "

# Row 14 - WrappedBTC
Set-TextValue "D14" "68.958.77"
$ws.Range("E14").Value = "  +0.54%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  -1.47%  "

# Row 16 - Avalanche
Set-TextValue "D16" "23.68"
$ws.Range("E16").Value = "  +0.03%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.826.98"
$ws.Range("E17").Value = "  +13.94%  "

# Row 18 - Chainlink
Set-TextValue "D18" "10.72"
$ws.Range("E18").Value = "  -2.65%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "336.76"
$ws.Range("E19").Value = "  -2.50%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.96"
$ws.Range("E20").Value = "  -3.17%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  -0.83%  "

# Row 22 - SuiNetwork
Set-TextValue "D22" "1.89"
$ws.Range("E22").Value = "  -0.33%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.11%  "

# Row 24 - Litecoin
Set-TextValue "D24" "66.83"
$ws.Range("E24").Value = "  -2.13%  "

# Row 25 - NEARProtocol
Set-TextValue "D25" "3.67"
$ws.Range("E25").Value = "  -2.70%  "

# Row 26 - Aptos
Set-TextValue "D26" "8.23"
$ws.Range("E26").Value = "  -0.31%  "

# Row 27 - PEPE
$ws.Range("D27").Value = "0.0₃0819"
$ws.Range("E27").Value = "  -2.86%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  -1.34%  "

# Rows 29/30 - swapped: Bittensor moves to 29, FirstDigitalUSD moves to 30
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D29" "433.95"
$ws.Range("E29").Value = "  -0.86%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  +0.00%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  -3.76%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -4.38%  "

# Row 33 - Monero
Set-TextValue "D33" "159.22"
$ws.Range("E33").Value = "  +1.20%  "

# Row 34 - WhiteBITCoin
$ws.Range("E34").Value = "  +0.06%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +0.25%  "

# Row 36 - USDe
$ws.Range("E36").Value = "  -0.06%  "

# Row 37 - EthereumClassic
Set-TextValue "D37" "17.81"
$ws.Range("E37").Value = "  -0.64%  "

# Row 38 - PolygonEcosystemToken
$ws.Range("E38").Value = "  -2.18%  "

# Row 39 - RenderToken
$ws.Range("E39").Value = "  -1.69%  "

# Row 40 - Stacks
Set-TextValue "D40" "1.47"
$ws.Range("E40").Value = "  -4.27%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  -2.88%  "

# Row 42 - dogwifhat
$ws.Range("E42").Value = "  -1.26%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  -0.76%  "

# Row 44 - Aave
Set-TextValue "D44" "131.38"
$ws.Range("E44").Value = "  -2.98%  "

# Rows 45/46 - swapped: Cronos moves to 45, ARBITRUM moves to 46
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D45" "0.0713"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D46" "0.485"
$ws.Range("E46").Value = "  -0.74%  "

# Row 47 - Mantle
Set-TextValue "D47" "0.564"
$ws.Range("E47").Value = "  -0.29%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -0.29%  "

# Row 49 - BitgetToken
$ws.Range("E49").Value = "  +0.14%  "

# Row 50 - Optimism
$ws.Range("E50").Value = "  -2.15%  "

# Row 51 - THORChain
$ws.Range("E51").Value = "  -7.73%  "
